# Scheduled market-data refresh: update cached currentAveragePrice /
# LevePrice / LeveProfit figures per leve row across all job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2864.6924
$ws.Range("I15").Value = 2864.6924
$ws.Range("K15").Value = 8594.0772
$ws.Range("M15").Value = -8425.0772

$ws.Range("H62").Value = 18523018
$ws.Range("I62").Value = 22226622
$ws.Range("K62").Value = 22226622
$ws.Range("M62").Value = -22225998

$ws.Range("H65").Value = 18523018
$ws.Range("I65").Value = 22226622
$ws.Range("K65").Value = 111133110
$ws.Range("M65").Value = -111129990

$ws.Range("H111").Value = 4603.222
$ws.Range("I111").Value = 6514.5
$ws.Range("J111").Value = 4057.1428
$ws.Range("K111").Value = 19543.5
$ws.Range("L111").Value = 12171.4284
$ws.Range("M111").Value = -16476.5
$ws.Range("N111").Value = -18305.4284

$ws.Range("H113").Value = 3651.5
$ws.Range("I113").Value = 3636.6667
$ws.Range("J113").Value = 3666.3333
$ws.Range("K113").Value = 3636.6667
$ws.Range("L113").Value = 3666.3333
$ws.Range("M113").Value = -382.6667000000002
$ws.Range("N113").Value = -10174.3333

$ws.Range("H116").Value = 2850.1177
$ws.Range("I116").Value = 2006.6666
$ws.Range("J116").Value = 3310.182
$ws.Range("K116").Value = 2006.6666
$ws.Range("L116").Value = 3310.182
$ws.Range("M116").Value = 1435.3334
$ws.Range("N116").Value = -10194.182

$ws.Range("H125").Value = 1981.55
$ws.Range("J125").Value = 1719.4706
$ws.Range("L125").Value = 15475.2354
$ws.Range("N125").Value = -20395.2354

$ws.Range("H138").Value = 1106.54
$ws.Range("I138").Value = 619.2174
$ws.Range("J138").Value = 1521.6666
$ws.Range("K138").Value = 1857.6522
$ws.Range("L138").Value = 4564.9998
$ws.Range("M138").Value = 3282.3478
$ws.Range("N138").Value = -14844.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4480.781
$ws.Range("I32").Value = 3962.9546
$ws.Range("J32").Value = 9363.143
$ws.Range("K32").Value = 3962.9546
$ws.Range("L32").Value = 9363.143
$ws.Range("M32").Value = -3675.9546
$ws.Range("N32").Value = -9937.143

$ws.Range("H45").Value = 1351.1666
$ws.Range("I45").Value = 1488.8889
$ws.Range("J45").Value = 938
$ws.Range("K45").Value = 1488.8889
$ws.Range("L45").Value = 938
$ws.Range("M45").Value = -1111.8889
$ws.Range("N45").Value = -1692

$ws.Range("H74").Value = 1626.8422
$ws.Range("I74").Value = 1126.4
$ws.Range("K74").Value = 1126.4
$ws.Range("M74").Value = -252.4000000000001

$ws.Range("H77").Value = 1626.8422
$ws.Range("I77").Value = 1126.4
$ws.Range("K77").Value = 5632
$ws.Range("M77").Value = -1264

$ws.Range("H132").Value = 1068.035
$ws.Range("I132").Value = 982.8570999999999
$ws.Range("K132").Value = 2948.5713
$ws.Range("M132").Value = -418.5712999999996

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 19231400
$ws.Range("I94").Value = 22727564
$ws.Range("K94").Value = 22727564
$ws.Range("M94").Value = -22727113

$ws.Range("H107").Value = 1647.5883
$ws.Range("I107").Value = 1331.2727
$ws.Range("J107").Value = 2227.5
$ws.Range("K107").Value = 1331.2727
$ws.Range("L107").Value = 2227.5
$ws.Range("M107").Value = 588.7273
$ws.Range("N107").Value = -6067.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2242.9678
$ws.Range("I31").Value = 2061.8333
$ws.Range("K31").Value = 2061.8333
$ws.Range("M31").Value = -1766.8333

$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()

$ws.Range("H34").Value = 2242.9678
$ws.Range("I34").Value = 2061.8333
$ws.Range("K34").Value = 2061.8333
$ws.Range("M34").Value = -1859.8333

$ws.Range("H58").Value = 688.86206
$ws.Range("I58").Value = 596.04
$ws.Range("K58").Value = 596.04
$ws.Range("M58").Value = -393.04

$ws.Range("H107").Value = 669.8
$ws.Range("I107").Value = 534.7
$ws.Range("J107").Value = 940
$ws.Range("K107").Value = 534.7
$ws.Range("L107").Value = 940
$ws.Range("M107").Value = 1385.3
$ws.Range("N107").Value = -4780

$ws.Range("H132").Value = 2985.8813
$ws.Range("I132").Value = 3200.75
$ws.Range("K132").Value = 9602.25
$ws.Range("M132").Value = -7072.25

$ws.Range("H134").Value = 10205341
$ws.Range("I134").Value = 1319.8292
$ws.Range("K134").Value = 3959.487599999999
$ws.Range("M134").Value = -1424.487599999999

$ws.Range("H136").Value = 688.86206
$ws.Range("I136").Value = 596.04
$ws.Range("K136").Value = 1788.12
$ws.Range("M136").Value = 761.8800000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 5359.3105
$ws.Range("J88").Value = 5697.037
$ws.Range("L88").Value = 17091.111
$ws.Range("N88").Value = -17947.111

$ws.Range("H91").Value = 5359.3105
$ws.Range("J91").Value = 5697.037
$ws.Range("L91").Value = 17091.111
$ws.Range("N91").Value = -20055.111

$ws.Range("H113").Value = 621.6905
$ws.Range("J113").Value = 653.3714
$ws.Range("L113").Value = 1960.1142
$ws.Range("N113").Value = -6300.1142

$ws.Range("H139").Value = 2116.2
$ws.Range("I139").Value = 2487.7368
$ws.Range("J139").Value = 1675
$ws.Range("K139").Value = 7463.2104
$ws.Range("L139").Value = 5025
$ws.Range("M139").Value = -2323.2104
$ws.Range("N139").Value = -15305

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H106").Value = 65000
$ws.Range("J106").Value = 65000
$ws.Range("L106").Value = 65000
$ws.Range("N106").Value = -67524

$ws.Range("H113").Value = 1724.6
$ws.Range("I113").Value = 1780.75
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1780.75
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 389.25
$ws.Range("N113").Value = -5840

$ws.Range("H122").Value = 1739.35
$ws.Range("I122").Value = 1488.1666
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 4464.4998
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -2014.4998
$ws.Range("N122").Value = -16900

$ws.Range("H124").Value = 45399
$ws.Range("J124").Value = 45399
$ws.Range("L124").Value = 45399
$ws.Range("N124").Value = -55219

$ws.Range("H132").Value = 1492.7222
$ws.Range("I132").Value = 1243.1333
$ws.Range("J132").Value = 2740.6667
$ws.Range("K132").Value = 3729.3999
$ws.Range("L132").Value = 8222.000100000001
$ws.Range("M132").Value = -1199.3999
$ws.Range("N132").Value = -13282.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2873.5
$ws.Range("I40").Value = 2569.7144
$ws.Range("K40").Value = 2569.7144
$ws.Range("M40").Value = -2433.7144

$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").ClearContents()

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()

$ws.Range("H132").Value = 23939.533
$ws.Range("I132").Value = 1461.2069
$ws.Range("J132").Value = 64681.5
$ws.Range("K132").Value = 4383.620699999999
$ws.Range("L132").Value = 194044.5
$ws.Range("M132").Value = -1853.620699999999
$ws.Range("N132").Value = -199104.5

$ws.Range("H136").Value = 3981.3428
$ws.Range("I136").Value = 4010.2058
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 12030.6174
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -9480.617400000001
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 635.5
$ws.Range("I107").Value = 530.6667
$ws.Range("J107").Value = 950
$ws.Range("K107").Value = 1592.0001
$ws.Range("L107").Value = 2850
$ws.Range("M107").Value = 327.9999
$ws.Range("N107").Value = -6690

$ws.Range("H124").Value = 65000
$ws.Range("J124").Value = 65000
$ws.Range("L124").Value = 65000
$ws.Range("N124").Value = -74820

$ws.Range("H132").Value = 1980.0952
$ws.Range("I132").Value = 2533.2
$ws.Range("J132").Value = 1477.2727
$ws.Range("K132").Value = 7599.599999999999
$ws.Range("L132").Value = 4431.8181
$ws.Range("M132").Value = -5069.599999999999
$ws.Range("N132").Value = -9491.8181

$ws.Range("H136").Value = 528.2381
$ws.Range("I136").Value = 502.05264
$ws.Range("J136").Value = 777
$ws.Range("K136").Value = 1506.15792
$ws.Range("L136").Value = 2331
$ws.Range("M136").Value = 1043.84208
$ws.Range("N136").Value = -7431
